# Sprint 42 test-case report update:
# - Fill in the "Day 9" (rows 50-53) summary counts that were previously blank.
# - Update the sheet's active selection / scroll position to reflect where the
#   user ended up working (C53, scrolled down to around row 41).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Day 9 summary values (previously empty cells) ---
$ws.Range("C51").Value = 6936
$ws.Range("C52").Value = 2182
$ws.Range("C53").Value = 2182

# --- View / selection state ---
# Scroll the window so row 41 is at the top-left, then select C53
# (mirrors the user scrolling down and clicking into the newly filled cell).
$win = $excel.ActiveWindow
$win.ScrollRow = 41
$win.ScrollColumn = 1
$ws.Range("C53").Select()
